$d = $word.ActiveDocument
$vt = [char]11   # <w:br/> manual line break inside a run

# ------------------------------------------------------------------
# 1) Delete whole paragraphs the revision drops outright. Deleting
#    from the bottom of the document upward keeps the 1-based
#    Paragraphs indices of earlier items valid while we go.
# ------------------------------------------------------------------

# Paragraphs 9 .. 18 (1-based): "תגמול אקראי..." through "סיכום:"
$blockStart = $d.Paragraphs.Item(9)
$blockEnd = $d.Paragraphs.Item(18)
$d.Range($blockStart.Range.Start, $blockEnd.Range.End).Delete()

# Heading3 "מה עשו החוקרים" (was #5)
$d.Paragraphs.Item(5).Range.Delete()

# Heading3 "למה זה חשוב" (was #3)
$d.Paragraphs.Item(3).Range.Delete()

# ------------------------------------------------------------------
# 2) Rewrite the text of the paragraphs that remain, via Find/Replace
#    scoped to each paragraph's own Range (keeps xml:space handling
#    identical to a plain Word find & replace -- no stray
#    xml:space="preserve" the way a raw Range.Text assignment adds).
#    After the deletes above the survivors are items 1..8.
# ------------------------------------------------------------------

# Paragraph 1: title / subtitle (2 runs split by an existing <w:br/>)
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute('המאמר היומי של יניב ומייק: 09.06.25', $true, $false, $false, $false, $false, $true, 1, $false, 'המאמר היומי של מייק: 07.06.25', 2) | Out-Null
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute('Spurious Rewards: Rethinking Training Signals in RLVR – Fast Overview', $true, $false, $false, $false, $false, $true, 1, $false, 'Rate-In: Information-Driven Adaptive Dropout Rates for Improved Inference-Time Uncertainty Estimation', 2) | Out-Null

# Paragraph 2: collapse the two old runs into one new sentence (wildcard
# spans the old <w:br/> between them)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Find.Execute('המסר המרכזי במשפט אחד*גם תגמולים אקראיים או שגויים יכולים להביא לשיפור דרמטי ביכולות פתרון בעיות מתמטיות – אבל רק אם המודל כבר "מכיר" את הדרך מהפרה-טריינינג.', $true, $false, $true, $false, $false, $true, 1, $false, 'היום אני סוקר מאמר מיוחד בכמה רבדים. הרובד הראשון אחד ממחבריו של מאמר זה הוא לא אחר אלא יאן לקון, אחד האבות של למידה עמוקה. הרובד השני מכיל את החוקר הישראלי הידוע (אך לא מספיק) רביד זיו שוורץ שהוא גם פרופסור באוניברסיטת ניו יורק. הרובד השלישי הוא נושא המאמר והוא שערוך אי ודאות עבור חיזוים של רשתות נוירונים - נושא שמאוד מעניין אותי אך לא מעט זמן לא סקרתי כזה.', 2) | Out-Null

# Paragraph 3 (was the RLVR blurb) -> new one-liner
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Find.Execute('למידה באמצעות חיזוקים עם תגמול ניתן לאימות (RL with Verifiable Rewards - RLVR) הפכה לשיטה מובילה לשפר יכולות חשיבה של מודלים גדולים. המאמר שואל שאלה פרובוקטיבית: האם אנחנו באמת צריכים תגמול מדויק? התשובה: לא תמיד.', $true, $false, $false, $false, $false, $true, 1, $false, 'איך ניתן לשערך אי הוודאות של החיזויים של רשת נוירונים? יש כמה משפחות של שיטות המוזכרות במאמר:', 2) | Out-Null

# Paragraph 4 (was "הם לקחו את המודל...") -> 5 lines joined by <w:br/>
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Find.Execute('הם לקחו את המודל Qwen-2.5-Math ואימנו אותו על סט שאלות מתמטיקה עם חמש גרסאות שונות של תגמולים:', $true, $false, $false, $false, $false, $true, 1, $false, ('רשתות נוירונים בייסיאניות מגדירות התפלגויות הסתברותיות על משקלי הרשת, מה שמאפשר למדל אי-ודאות דרך ההתפלגות הפוסטריורית. עם זאת, הן כבדות חישובית וקשה להרחיב אותן.' + $vt + 'שיטות אנסמבל: מאמנות מספר מודלים ומאגדות את התחזיות שלהם. מסוגלות למדל גם אי ודאות אפיסטמית וגם אליאטורית, אך דורשות משאבים חישוביים רבים.' + $vt + 'אוגמנטציה של דאטה בזמן טסט (Test-Time Augmentation): מוסיפות שיבושים לקלט (כמו סיבוב או טשטוש) כדי להעריך את התפלגות התחזיות. יעיל בעיקר כשיש ידע מוקדם על מבנה הנתונים.' + $vt + 'הזרקת רעש למודל: מוסיפים רעש נשלט (למשל גאוסי) למשקלים או לפעולות כדי לבחון רגישות מעבר לשינויים בקלט.' + $vt + 'שיטות מונטה קרלו (MC): משתמשות בדגימות אקראיות כדי לאמוד אי ודאות. למשל, MC Dropout מפעיל דרופאוט(dropout) גם בזמן טסט כדי לדגום את מרחב משקלי הרשת. יש לא מעט שיטות נוספות מבוססת MC לשערוך אי ודאות ברשתות.'), 2) | Out-Null

# Paragraph 5 (was "תגמול אמיתי...") -> intro line, blank <w:br/><w:br/>, then 3 lines
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Find.Execute('תגמול אמיתי:  מודל מקבל נקודה רק אם התשובה נכונה.', $true, $false, $false, $false, $false, $true, 1, $false, ('אבל איך ניתן לשערך את הוודאות? אחת הדרכים היא להשתמש בגישות מתורת המידע (information theory) לניתוח של זרימה המידע בתוך הרשת ומידת ״פגיעתה״ מהשיטות המוזכרות מעלה (למשל MC Dropout). בגדול מאוד ככל שזרימת המידע נפגעת יותר - אי הוודאות של החיזויים עולה. שיטות מתורת המידע די נפוצות במחקר של רשתות עמוקות למשל:' + $vt + $vt + 'עקרון צוואר הבקבוק המידעי (של נפתלי תשבי): מציע ששכבות ברשת נוירונים שואפות לדחוס את המידע מהקלט תוך שמירה על המידע הרלוונטי לפלט. משמש לניתוח דינמיקת הלמידה והכללה של המודל.' + $vt + 'ניתוח מידע הדדי (Mutual Information): הערכת המידע ההדדי בין הקלט, השכבות הפנימיות והפלט מסייעת להבין כיצד מידע זורם ומעובד ברשת. זה הטכניקה שהמחברים משתמשים בה במאמר' + $vt + 'טכניקות רגולריזציה אינפורמטיביות: שיטות כמו information dropout שולטות בזרימת המידע במהלך האימון כדי לשפר חוסן והכללה של המודל.'), 2) | Out-Null

# Paragraph 6 (was two runs: "תגמול לפי הצבעת רוב..." + "תגמול פורמטי...")
# -> a single new paragraph (trailing space -> xml:space preserved automatically).
# The wildcard prefix/suffix avoid the \boxed{} special wildcard characters in
# the middle of the old text -- '*' alone is enough to span the old <w:br/>.
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Find.Execute('תגמול לפי הצבעת רוב: המודל מייצר 64 תשובות, ומתגמל את התשובה השכיחה.*בלי קשר לנכונות.', $true, $false, $true, $false, $false, $true, 1, $false, 'אוקיי, אז המאמר מציע שיטה מבוססת מידע הדדי המשכללת MC dropout. במקום להשתמש ב dropout rate קבוע לכל השכבות (כלומר מה אחוז הנוירונים המחוסלים בשכבה) המחברים מציע לקבוע אותה (dropout rate) בתתלות במידת פגיעתה בזרימת המידע בשכבה. המטרה כאן היא לעשות את אובדן המידע בכל שכבה פחות או יותר קבוע. אם אובדן המידע הדדי גבוה(מקבוע אפסילון) מדי מקטינים dropout rate ואם זה נמוך מדי מגדילים אותו. ', 2) | Out-Null

# Paragraph 7 (was the closing "המאמר מראה..." summary) -> new closing remark
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Find.Execute('המאמר מראה שלעיתים קרובות אימון RL לא מלמד כישורים חדשים, אלא מחלץ כישורים חבויים שהמודל כבר פיתח בפרה-טריינינג. לא תמיד צריך תגמול מדויק – אם המודל כבר "מכיר" את הדרך, מספיק לאותת לו לחזור אליה. עם זאת, זה לא נכון לכל מודל – יש כאלה שדורשים הנחיה מדויקת כדי להשתפר.', $true, $false, $false, $false, $false, $true, 1, $false, 'ד״א פגיעה בזרימת המידע בשכבה מחושבת דרך חישוב של המידע הדדי בין אקטיבציות של הקלט בשכבה לבין אלו של פלט השכבה. מתברר שזה די לא טריוויאלי והמאמר דן בהרחבה איך ניתן לעשות זאת.', 2) | Out-Null

# Paragraph 8 (the arxiv URL) is untouched.

Write-Output 'edit complete'
